# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q2" and "总计".
# We insert a new "2022-Q1" sheet (per-fund holding detail) right after
# "2021-Q2" -- cloned from "总计" so it picks up the same header/row-index
# cell styling -- and update the "总计" summary sheet so it also lists the
# new 2022-Q1 quarter (while keeping the prior 2021-Q2 row).

$wb = $excel.ActiveWorkbook
$q2Sheet = $wb.Worksheets.Item("2021-Q2")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Insert the new "2022-Q1" sheet right after "2021-Q2" ---------------
# Cloning "总计" (rather than Worksheets.Add()) means the new sheet starts
# out with the same bold/bordered header + row-index styling already baked
# in, instead of plain default formatting.
$totalSheet.Copy($null, $q2Sheet)

# NOTE: inserting a sheet shifts everyone after it, and sheet handles here
# track *position*, not identity -- re-resolve both sheets by name now that
# the sheet count/order has changed, instead of trusting old handles.
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item("总计")

# Extend that header styling across the new E:H columns.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "090019"
$newSheet.Range("C2").Value = "大成景恒混合A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.31"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.51"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.93"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0446"
$newSheet.Range("H2").Value = 7

# Row 3 is brand new -- clone A2's row-index styling onto A3 first.
$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false
$newSheet.Range("A3").Value = 1

$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "006038"
$newSheet.Range("C3").Value = "大成景恒混合C"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.92"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "93.51"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "1.93"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0178"
$newSheet.Range("H3").Value = 7

# --- 2. Update the "总计" summary sheet to also include 2022-Q1 ------------
# Row 3 is brand new -- clone A2's row-index styling onto A3 first, then
# move the existing 2021-Q2 figures down into it.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Application.CutCopyMode = $false

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.01

# Row 2 becomes the new 2022-Q1 totals.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.06

$q2Sheet.Activate()
